# tests/fixtures/templates/tags/td/td.xlsx
# fix(MSE-1120): Fix <td> child recursion.
#
# The fixture was regenerated after fixing how the <td> HTML tag's children
# are walked/recursed when building the worksheet. Net effect on this sample
# sheet:
#   - the sheet got a fresh generated name
#   - the second <td> ("col 1" in A2), which previously fell through to the
#     workbook's default/Normal formatting, now correctly inherits the
#     sheet's "bold" (Arial, 10pt, bold, no special colour) text formatting
#     with a text ("@") number format - matching how the first row's <td>
#     cells (A1/B1) are handled
#   - the active selection on open is back at A1 instead of A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated template id for this sheet.
$ws.Name = "7e862dbb"

# A2 ("col 1") now picks up the bold Arial 10pt formatting (same font as the
# header cells, just without the red colour) plus a text number format,
# instead of staying on the default/Normal style.
$ws.Range("A2").Style = "bold"
$ws.Range("A2").NumberFormat = "@"

# Selection reverts to A1 (previously A2).
$null = $ws.Range("A1").Select()
